$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the paragraph that holds the "${picture6}" placeholder run. We walk
# the Paragraphs collection (instead of trusting a hard-coded index) so the
# script keeps working even if earlier content in the document shifts the
# paragraph numbering around.
# ---------------------------------------------------------------------------
$picture6Para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*`${picture6}*") {
        $picture6Para = $p
        break
    }
}

if ($picture6Para -ne $null) {
    # Delete the paragraph mark that ends the "${picture6}" paragraph. This
    # merges it with the following paragraph (the "${picture7}" one); the
    # merged paragraph keeps the *second* paragraph's formatting, which is
    # exactly what turns the first paragraph's rFonts hint from "eastAsia"
    # into "default".
    $mark1 = $d.Range($picture6Para.Range.End - 1, $picture6Para.Range.End)
    $mark1.Delete()

    # Re-fetch the (now merged) paragraph and delete its new trailing mark
    # too. That merges in the third paragraph, which only contained the
    # "_GoBack" bookmark - pulling the bookmark up into the same paragraph
    # as the "${picture6}"/"${picture7}" runs.
    $mergedPara = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*`${picture6}*") {
            $mergedPara = $p
            break
        }
    }
    $mark2 = $d.Range($mergedPara.Range.End - 1, $mergedPara.Range.End)
    $mark2.Delete()
}

# Remove the now-redundant "${picture7}" run text, leaving the bookmark
# alone in the merged paragraph together with the "${picture6}" run.
$d.Content.Find.Execute("`${picture7}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 1)
